# WorkoutProgram.xlsx — rework the exercise list for DAY 1 - DAY 5,
# and remove DAY 6 entirely (per "Add files via upload" re-upload of the
# workout plan with a new exercise rotation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the whole exercise block (rows 5-38): old content had DAY 1..DAY 6
# each with a variable number of exercise rows; the new content only goes
# through DAY 5 and the row layout per day is different, so it's cleanest to
# wipe the block and rewrite it fresh rather than trying to patch individual
# rows in place.
$ws.Range("A5:M38").Clear()

# ---- DAY 1 ----
$ws.Range("A5").Value = "DAY 1"
$ws.Range("A6").Value = "Exercises"
$ws.Range("B6").Value = "Sets"
$ws.Range("C6").Value = "Reps"
$ws.Range("A7").Value = " Smith Machine Shoulder Press"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 8

# ---- DAY 2 ----
$ws.Range("A9").Value = "DAY 2"
$ws.Range("A10").Value = "Exercises"
$ws.Range("B10").Value = "Sets"
$ws.Range("C10").Value = "Reps"
$ws.Range("A11").Value = "Cable EZ Bar Curl"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 8
$ws.Range("A12").Value = "Shrug Machine"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 8

# ---- DAY 3 ----
$ws.Range("A14").Value = "DAY 3"
$ws.Range("A15").Value = "Exercises"
$ws.Range("B15").Value = "Sets"
$ws.Range("C15").Value = "Reps"
$ws.Range("A16").Value = " Deficit Smith Machine Calves"
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = 8
$ws.Range("A17").Value = "Machine Hip Thrust"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 8
$ws.Range("A18").Value = "Leg Extensions"
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 8

# ---- DAY 3 (second block, same label as the source workbook) ----
$ws.Range("A20").Value = "DAY 3"
$ws.Range("A21").Value = "Exercises"
$ws.Range("B21").Value = "Sets"
$ws.Range("C21").Value = "Reps"
$ws.Range("A22").Value = "Shoulder Press Machine"
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 8

# ---- DAY 4 ----
$ws.Range("A24").Value = "DAY 4"
$ws.Range("A25").Value = "Exercises"
$ws.Range("B25").Value = "Sets"
$ws.Range("C25").Value = "Reps"
$ws.Range("A26").Value = "Cable EZ Bar Curl"
$ws.Range("B26").Value = 4
$ws.Range("C26").Value = 8
$ws.Range("A27").Value = "Shrug Machine"
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 8

# ---- DAY 5 ----
$ws.Range("A29").Value = "DAY 5"
$ws.Range("A30").Value = "Exercises"
$ws.Range("B30").Value = "Sets"
$ws.Range("C30").Value = "Reps"
$ws.Range("A31").Value = "Calf Machine"
$ws.Range("B31").Value = 4
$ws.Range("C31").Value = 8
$ws.Range("A32").Value = " Cable Pull through"
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = 8
$ws.Range("A33").Value = "Leg Extensions"
$ws.Range("B33").Value = 4
$ws.Range("C33").Value = 8
